$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" column header in H1, reusing the existing header style
# (same bold/centered/bordered look as the other header cells).
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# New "Save" column values (0/1 flags) for rows 2-12
$values = @(1, 0, 0, 1, 1, 1, 0, 0, 1, 0, 0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
